# Regenerate orders with updated distance/size codes.
# Distance codes: D51 -> D55, D64 -> D69, D80 -> D86
# Size code:      S30 -> S31
# These substrings appear embedded inside larger text tokens
# (e.g. "Face14_D51_S25", "Face08_D51_S30_l.png", "D51") across the
# used range of the sheet, so every text cell is scanned and rewritten
# in place with the substitutions applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val
            $newVal = $newVal -replace 'D51', 'D55'
            $newVal = $newVal -replace 'D64', 'D69'
            $newVal = $newVal -replace 'D80', 'D86'
            $newVal = $newVal -replace 'S30', 'S31'

            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
